$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.292101740837097
$ws.Range("B1").Value = 2.354964256286621
$ws.Range("C1").Value = 3.057630300521851
$ws.Range("D1").Value = 3.519988536834717
$ws.Range("E1").Value = 1.335824728012085
